$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Probabilities")

# Row 2
$ws.Cells.Item(2, 1).Value = 1369
$ws.Cells.Item(2, 2).Value = "2025-11-26T12:30:00"
$ws.Cells.Item(2, 3).Value = "Адмирал"
$ws.Cells.Item(2, 4).Value = "Амур"
$ws.Cells.Item(2, 5).Value = 897818
$ws.Cells.Item(2, 6).Value = "https://text.khl.ru/text/897818.html"
$ws.Cells.Item(2, 7).Value = 1.91778
$ws.Cells.Item(2, 8).Value = 2.626433
$ws.Cells.Item(2, 9).Value = 1.554363
$ws.Cells.Item(2, 10).Value = 4.8
$ws.Cells.Item(2, 11).Value = 3.35889
$ws.Cells.Item(2, 12).Value = 2.090398
$ws.Cells.Item(2, 13).Value = 4.544213
$ws.Cells.Item(2, 14).Value = 29.99177
$ws.Cells.Item(2, 15).Value = 28.61054
$ws.Cells.Item(2, 16).Value = 58.60231
$ws.Cells.Item(2, 17).Value = -0.127799
$ws.Cells.Item(2, 18).Value = 0.050246
$ws.Cells.Item(2, 19).Value = 0.625932
$ws.Cells.Item(2, 20).Value = 0.153208
$ws.Cells.Item(2, 21).Value = 0.220113
$ws.Cells.Item(2, 22).Value = 0.207513
$ws.Cells.Item(2, 23).Value = 0.79174
$ws.Cells.Item(2, 24).Value = 0.365474
$ws.Cells.Item(2, 25).Value = 0.633778
$ws.Cells.Item(2, 26).Value = 0.53763
$ws.Cells.Item(2, 27).Value = 0.461622
$ws.Cells.Item(2, 28).Value = 0.693985
$ws.Cells.Item(2, 29).Value = 0.305268
$ws.Cells.Item(2, 30).Value = 0.815702
$ws.Cells.Item(2, 31).Value = 0.18355
$ws.Cells.Item(2, 32).Value = 0.848425
$ws.Cells.Item(2, 33).Value = 0.151575
$ws.Cells.Item(2, 34).Value = 0.652263
$ws.Cells.Item(2, 35).Value = 0.347737
$ws.Cells.Item(2, 36).Value = 0.61791
$ws.Cells.Item(2, 37).Value = 0.38209
$ws.Cells.Item(2, 38).Value = 0.347775
$ws.Cells.Item(2, 39).Value = 0.652225
$ws.Cells.Item(2, 40).Value = 0.887905
$ws.Cells.Item(2, 41).Value = 0.548087

# Row 3
$ws.Cells.Item(3, 1).Value = 1369
$ws.Cells.Item(3, 2).Value = "2025-11-26T17:00:00"
$ws.Cells.Item(3, 3).Value = "Трактор"
$ws.Cells.Item(3, 4).Value = "Драконы"
$ws.Cells.Item(3, 5).Value = 897816
$ws.Cells.Item(3, 6).Value = "https://text.khl.ru/text/897816.html"
$ws.Cells.Item(3, 7).Value = 4.435928
$ws.Cells.Item(3, 8).Value = 3.568607
$ws.Cells.Item(3, 9).Value = 5.035714
$ws.Cells.Item(3, 10).Value = 2.320284
$ws.Cells.Item(3, 11).Value = 3.378106
$ws.Cells.Item(3, 12).Value = 4.302161
$ws.Cells.Item(3, 13).Value = 8.004535
$ws.Cells.Item(3, 14).Value = 39.059641
$ws.Cells.Item(3, 15).Value = 31.713169
$ws.Cells.Item(3, 16).Value = 70.77281
$ws.Cells.Item(3, 17).Value = 0.140026
$ws.Cells.Item(3, 18).Value = 0.131141
$ws.Cells.Item(3, 19).Value = 0.300991
$ws.Cells.Item(3, 20).Value = 0.139072
$ws.Cells.Item(3, 21).Value = 0.554329
$ws.Cells.Item(3, 22).Value = 0.052503
$ws.Cells.Item(3, 23).Value = 0.941889
$ws.Cells.Item(3, 24).Value = 0.11946
$ws.Cells.Item(3, 25).Value = 0.874932
$ws.Cells.Item(3, 26).Value = 0.22231
$ws.Cells.Item(3, 27).Value = 0.772082
$ws.Cells.Item(3, 28).Value = 0.353962
$ws.Cells.Item(3, 29).Value = 0.64043
$ws.Cells.Item(3, 30).Value = 0.498408
$ws.Cells.Item(3, 31).Value = 0.495984
$ws.Cells.Item(3, 32).Value = 0.850654
$ws.Cells.Item(3, 33).Value = 0.149346
$ws.Cells.Item(3, 34).Value = 0.656018
$ws.Cells.Item(3, 35).Value = 0.343982
$ws.Cells.Item(3, 36).Value = 0.928213
$ws.Cells.Item(3, 37).Value = 0.071787
$ws.Cells.Item(3, 38).Value = 0.802916
$ws.Cells.Item(3, 39).Value = 0.197084
$ws.Cells.Item(3, 40).Value = 0.586316
$ws.Cells.Item(3, 41).Value = 0.80824

# Row 4
$ws.Cells.Item(4, 1).Value = 1369
$ws.Cells.Item(4, 2).Value = "2025-11-26T19:00:00"
$ws.Cells.Item(4, 3).Value = "Северсталь"
$ws.Cells.Item(4, 4).Value = "СКА"
$ws.Cells.Item(4, 5).Value = 897817
$ws.Cells.Item(4, 6).Value = "https://text.khl.ru/text/897817.html"
$ws.Cells.Item(4, 7).Value = 1.464286
$ws.Cells.Item(4, 8).Value = 2.390395
$ws.Cells.Item(4, 9).Value = 1.107143
$ws.Cells.Item(4, 10).Value = 1.410664
$ws.Cells.Item(4, 11).Value = 1.437475
$ws.Cells.Item(4, 12).Value = 1.748769
$ws.Cells.Item(4, 13).Value = 3.854681
$ws.Cells.Item(4, 14).Value = 24.993415
$ws.Cells.Item(4, 15).Value = 27.964901
$ws.Cells.Item(4, 16).Value = 52.958316
$ws.Cells.Item(4, 17).Value = -0.2
$ws.Cells.Item(4, 18).Value = 0.009782
$ws.Cells.Item(4, 19).Value = 0.317655
$ws.Cells.Item(4, 20).Value = 0.231914
$ws.Cells.Item(4, 21).Value = 0.450429
$ws.Cells.Item(4, 22).Value = 0.605583
$ws.Cells.Item(4, 23).Value = 0.394414
$ws.Cells.Item(4, 24).Value = 0.783058
$ws.Cells.Item(4, 25).Value = 0.216939
$ws.Cells.Item(4, 26).Value = 0.896154
$ws.Cells.Item(4, 27).Value = 0.103844
$ws.Cells.Item(4, 28).Value = 0.956212
$ws.Cells.Item(4, 29).Value = 0.043785
$ws.Cells.Item(4, 30).Value = 0.983549
$ws.Cells.Item(4, 31).Value = 0.016448
$ws.Cells.Item(4, 32).Value = 0.421034
$ws.Cells.Item(4, 33).Value = 0.578966
$ws.Cells.Item(4, 34).Value = 0.17563
$ws.Cells.Item(4, 35).Value = 0.82437
$ws.Cells.Item(4, 36).Value = 0.521747
$ws.Cells.Item(4, 37).Value = 0.478253
$ws.Cells.Item(4, 38).Value = 0.255703
$ws.Cells.Item(4, 39).Value = 0.744297
$ws.Cells.Item(4, 40).Value = 0.7598
$ws.Cells.Item(4, 41).Value = 0.855151

# Row 5
$ws.Cells.Item(5, 1).Value = 1369
$ws.Cells.Item(5, 2).Value = "2025-11-26T19:30:00"
$ws.Cells.Item(5, 3).Value = "Динамо М"
$ws.Cells.Item(5, 4).Value = "Локомотив"
$ws.Cells.Item(5, 5).Value = 897815
$ws.Cells.Item(5, 6).Value = "https://text.khl.ru/text/897815.html"
$ws.Cells.Item(5, 7).Value = 2.014963
$ws.Cells.Item(5, 8).Value = 2.067992
$ws.Cells.Item(5, 9).Value = 2.713803
$ws.Cells.Item(5, 10).Value = 1.78325
$ws.Cells.Item(5, 11).Value = 1.899107
$ws.Cells.Item(5, 12).Value = 2.390898
$ws.Cells.Item(5, 13).Value = 4.082955
$ws.Cells.Item(5, 14).Value = 24.33595
$ws.Cells.Item(5, 15).Value = 26.124622
$ws.Cells.Item(5, 16).Value = 50.460572
$ws.Cells.Item(5, 17).Value = -0.110672
$ws.Cells.Item(5, 18).Value = -0.124268
$ws.Cells.Item(5, 19).Value = 0.311642
$ws.Cells.Item(5, 20).Value = 0.194412
$ws.Cells.Item(5, 21).Value = 0.493899
$ws.Cells.Item(5, 22).Value = 0.378954
$ws.Cells.Item(5, 23).Value = 0.620999
$ws.Cells.Item(5, 24).Value = 0.572371
$ws.Cells.Item(5, 25).Value = 0.427582
$ws.Cells.Item(5, 26).Value = 0.738323
$ws.Cells.Item(5, 27).Value = 0.261631
$ws.Cells.Item(5, 28).Value = 0.856978
$ws.Cells.Item(5, 29).Value = 0.142975
$ws.Cells.Item(5, 30).Value = 0.929697
$ws.Cells.Item(5, 31).Value = 0.070256
$ws.Cells.Item(5, 32).Value = 0.565997
$ws.Cells.Item(5, 33).Value = 0.434003
$ws.Cells.Item(5, 34).Value = 0.296038
$ws.Cells.Item(5, 35).Value = 0.703962
$ws.Cells.Item(5, 36).Value = 0.689572
$ws.Cells.Item(5, 37).Value = 0.310428
$ws.Cells.Item(5, 38).Value = 0.427911
$ws.Cells.Item(5, 39).Value = 0.572089
$ws.Cells.Item(5, 40).Value = 0.696463
$ws.Cells.Item(5, 41).Value = 0.839553

